# Implement type designators for technique-specific sample preparation details.
#
# Adds the SamplePreparation "mixin" columns (preparation_type, sample_id,
# preparation_date, operator_id, protocol_description, id, title) to each of
# the technique-specific preparation sheets - CryoEMPreparation,
# XRayPreparation, SAXSPreparation - inserting them just before the existing
# trailing "description" column, and removes the now-redundant preparation_type
# drop-down validation from the base SamplePreparation sheet.

$wb = $excel.ActiveWorkbook

# xlShiftToRight
$ShiftToRight = -4161

$newCols = @(
    "preparation_type",
    "sample_id",
    "preparation_date",
    "operator_id",
    "protocol_description",
    "id",
    "title"
)

# --- CryoEMPreparation: A1:J1 -> A1:Q1 (insert before J1 "description") ---
$ws = $wb.Worksheets.Item("CryoEMPreparation")
$ws.Range("J1:P1").Insert($ShiftToRight)
$ws.Range("J1").Value = $newCols[0]
$ws.Range("K1").Value = $newCols[1]
$ws.Range("L1").Value = $newCols[2]
$ws.Range("M1").Value = $newCols[3]
$ws.Range("N1").Value = $newCols[4]
$ws.Range("O1").Value = $newCols[5]
$ws.Range("P1").Value = $newCols[6]

# --- XRayPreparation: A1:H1 -> A1:O1 (insert before H1 "description") ---
$ws = $wb.Worksheets.Item("XRayPreparation")
$ws.Range("H1:N1").Insert($ShiftToRight)
$ws.Range("H1").Value = $newCols[0]
$ws.Range("I1").Value = $newCols[1]
$ws.Range("J1").Value = $newCols[2]
$ws.Range("K1").Value = $newCols[3]
$ws.Range("L1").Value = $newCols[4]
$ws.Range("M1").Value = $newCols[5]
$ws.Range("N1").Value = $newCols[6]

# --- SAXSPreparation: A1:F1 -> A1:M1 (insert before F1 "description") ---
$ws = $wb.Worksheets.Item("SAXSPreparation")
$ws.Range("F1:L1").Insert($ShiftToRight)
$ws.Range("F1").Value = $newCols[0]
$ws.Range("G1").Value = $newCols[1]
$ws.Range("H1").Value = $newCols[2]
$ws.Range("I1").Value = $newCols[3]
$ws.Range("J1").Value = $newCols[4]
$ws.Range("K1").Value = $newCols[5]
$ws.Range("L1").Value = $newCols[6]

# --- SamplePreparation: drop the preparation_type dropdown validation ---
# (technique-specific subclasses now carry their own designated type, so the
#  generic enumerated list on the base sheet is removed)
$ws = $wb.Worksheets.Item("SamplePreparation")
$ws.Range("A2:A1048576").Validation.Delete()
